# Insert a new data row at row 118 (weekly price update), pushing all
# subsequent rows (old 118..226) down by one to (119..227).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("118").Insert()

$ws.Range("A118").Value = 8
$ws.Range("B118").Value = "Terminal La Palmera de La Serena"
$ws.Range("C118").Value = "Coquimbo"
$ws.Range("D118").Value = 44512
$ws.Range("E118").Value = 4
$ws.Range("F118").Value = 100114013
$ws.Range("G118").Value = "Zanahoria"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 800
$ws.Range("K118").Value = 6500
$ws.Range("L118").Value = 7000
$ws.Range("M118").Value = 6750
$ws.Range("N118").Value = "$/saco 20 kilos"
$ws.Range("O118").Value = "Provincia del Elquí"
$ws.Range("P118").Value = 338
$ws.Range("Q118").Value = 20
$ws.Range("R118").Value = "Hortaliza"
